$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 44874
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 7500
$ws.Range("O2").Value = 8000
$ws.Range("P2").Value = 7750
$ws.Range("S2").Value = 7750

# Row 3 updates
$ws.Range("D3").Value = 44881
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 11250
$ws.Range("O3").Value = 11250
$ws.Range("P3").Value = 11250
$ws.Range("S3").Value = 11250
